$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.442.45'
$ws.Range("E2").Value = '  -1.82%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.840.68'
$ws.Range("E3").Value = '  -2.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '260.84'
$ws.Range("E5").Value = '  -6.52%  '
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5218'
$ws.Range("E7").Value = '  -1.61%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3242'
$ws.Range("E8").Value = '  -6.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06810'
$ws.Range("E9").Value = '  -2.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.64'
$ws.Range("E10").Value = '  -7.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7666'
$ws.Range("E11").Value = '  -5.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07695'
$ws.Range("E12").Value = '  -1.71%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.836.68'
$ws.Range("E13").Value = '  -2.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.54'
$ws.Range("E14").Value = '  -2.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.029'
$ws.Range("E15").Value = '  -3.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.93'
$ws.Range("E17").Value = '  -4.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007953'
$ws.Range("E19").Value = '  -1.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.496.65'
$ws.Range("E20").Value = '  -1.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.076.06'
$ws.Range("E21").Value = '  -1.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.566'
$ws.Range("E22").Value = '  -4.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.944'
$ws.Range("E24").Value = '  -4.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.61'
$ws.Range("E25").Value = '  -1.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.219'
$ws.Range("E26").Value = '  -6.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.656'
$ws.Range("E27").Value = '  -0.47%  '
$ws.Range("E28").Value = '  -2.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.43'
$ws.Range("E29").Value = '  -2.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.164'
$ws.Range("E30").Value = '  -4.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.134'
$ws.Range("E31").Value = '  -4.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08735'
$ws.Range("E32").Value = '  -1.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04801'
$ws.Range("E33").Value = '  -3.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.121'
$ws.Range("E34").Value = '  -5.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.847'
$ws.Range("E35").Value = '  -1.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6992'
$ws.Range("E36").Value = '  -5.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.063'
$ws.Range("E37").Value = '  -7.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01761'
$ws.Range("E38").Value = '  -5.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.190'
$ws.Range("E39").Value = '  -8.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4833'
$ws.Range("E40").Value = '  -6.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '111.31'
$ws.Range("E41").Value = '  -4.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8904'
$ws.Range("E42").Value = '  -7.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.077'
$ws.Range("E43").Value = '  -2.02%  '
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.670'
$ws.Range("E45").Value = '  -5.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4137'
$ws.Range("E46").Value = '  -8.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05859'
$ws.Range("E47").Value = '  -1.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.019'
$ws.Range("E48").Value = '  -3.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.01'
$ws.Range("E49").Value = '  -3.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1221'
$ws.Range("E50").Value = '  -9.56%  '
$ws.Range("E51").Value = '  -0.43%  '
